# chore: update Sheets via scheduled runner
# Refresh cached market-board derived figures (currentAveragePrice*,
# LevePrice*/LeveProfit* columns H:N) for the rows whose underlying
# item prices changed since the last run.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 562.65717
$ws.Range("I80").Value = 551.75
$ws.Range("K80").Value = 1655.25
$ws.Range("M80").Value = -657.25
$ws.Range("H83").Value = 562.65717
$ws.Range("I83").Value = 551.75
$ws.Range("K83").Value = 4965.75
$ws.Range("M83").Value = 26.25
$ws.Range("H112").Value = 14286876
$ws.Range("J112").Value = 16130231
$ws.Range("L112").Value = 48390693
$ws.Range("N112").Value = -48392909
$ws.Range("H138").Value = 4870.906
$ws.Range("I138").Value = 996.7273
$ws.Range("J138").Value = 6223.7935
$ws.Range("K138").Value = 2990.1819
$ws.Range("L138").Value = 18671.3805
$ws.Range("M138").Value = 2149.8181
$ws.Range("N138").Value = -28951.3805

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5186.353
$ws.Range("I32").Value = 4011.1714
$ws.Range("J32").Value = 7757.0625
$ws.Range("K32").Value = 4011.1714
$ws.Range("L32").Value = 7757.0625
$ws.Range("M32").Value = -3724.1714
$ws.Range("N32").Value = -8331.0625
$ws.Range("H61").Value = 20204.834
$ws.Range("I61").Value = 23725.8
$ws.Range("J61").Value = 2600
$ws.Range("K61").Value = 23725.8
$ws.Range("L61").Value = 2600
$ws.Range("M61").Value = -23513.8
$ws.Range("N61").Value = -3024
$ws.Range("H123").Value = 33476.332
$ws.Range("J123").Value = 33476.332
$ws.Range("L123").Value = 33476.332
$ws.Range("N123").Value = -43276.332
$ws.Range("H132").Value = 5874.8857
$ws.Range("I132").Value = 1538.2759
$ws.Range("J132").Value = 26835.166
$ws.Range("K132").Value = 4614.8277
$ws.Range("L132").Value = 80505.49800000001
$ws.Range("M132").Value = -2084.8277
$ws.Range("N132").Value = -85565.49800000001
$ws.Range("H134").Value = 48429
$ws.Range("J134").Value = 48429
$ws.Range("L134").Value = 48429
$ws.Range("N134").Value = -58569
$ws.Range("H136").Value = 20204.834
$ws.Range("I136").Value = 23725.8
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 71177.39999999999
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -68627.39999999999
$ws.Range("N136").Value = -12900

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9019.588
$ws.Range("I134").Value = 17719
$ws.Range("J134").Value = 2930
$ws.Range("K134").Value = 53157
$ws.Range("L134").Value = 8790
$ws.Range("M134").Value = -50622
$ws.Range("N134").Value = -13860

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5945.4814
$ws.Range("I31").Value = 1786.7142
$ws.Range("J31").Value = 10424.154
$ws.Range("K31").Value = 1786.7142
$ws.Range("L31").Value = 10424.154
$ws.Range("M31").Value = -1491.7142
$ws.Range("N31").Value = -11014.154
$ws.Range("H34").Value = 5945.4814
$ws.Range("I34").Value = 1786.7142
$ws.Range("J34").Value = 10424.154
$ws.Range("K34").Value = 1786.7142
$ws.Range("L34").Value = 10424.154
$ws.Range("M34").Value = -1584.7142
$ws.Range("N34").Value = -10828.154
$ws.Range("H53").Value = 40000
$ws.Range("J53").Value = 40000
$ws.Range("L53").Value = 40000
$ws.Range("N53").Value = -41214
$ws.Range("H62").Value = 8269.444
$ws.Range("I62").Value = 7202.7144
$ws.Range("J62").Value = 12003
$ws.Range("K62").Value = 7202.7144
$ws.Range("L62").Value = 12003
$ws.Range("M62").Value = -6578.7144
$ws.Range("N62").Value = -13251
$ws.Range("H65").Value = 8269.444
$ws.Range("I65").Value = 7202.7144
$ws.Range("J65").Value = 12003
$ws.Range("K65").Value = 36013.572
$ws.Range("L65").Value = 60015
$ws.Range("M65").Value = -32893.572
$ws.Range("N65").Value = -66255
$ws.Range("H132").Value = 3203.8
$ws.Range("I132").Value = 2634.5
$ws.Range("J132").Value = 4342.4
$ws.Range("K132").Value = 7903.5
$ws.Range("L132").Value = 13027.2
$ws.Range("M132").Value = -5373.5
$ws.Range("N132").Value = -18087.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 356.2353
$ws.Range("I107").Value = 241.85715
$ws.Range("J107").Value = 436.3
$ws.Range("K107").Value = 725.5714499999999
$ws.Range("L107").Value = 1308.9
$ws.Range("M107").Value = 1194.42855
$ws.Range("N107").Value = -5148.9
$ws.Range("H129").Value = 1556.9
$ws.Range("J129").Value = 3533.3333
$ws.Range("L129").Value = 10599.9999
$ws.Range("N129").Value = -20599.9999
$ws.Range("H131").Value = 29744064
$ws.Range("I131").Value = 7178992.5
$ws.Range("J131").Value = 43479324
$ws.Range("K131").Value = 21536977.5
$ws.Range("L131").Value = 130437972
$ws.Range("M131").Value = -21531937.5
$ws.Range("N131").Value = -130448052

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7962.4
$ws.Range("I102").Value = 1906
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 1906
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -284
$ws.Range("N102").Value = -15244
$ws.Range("H107").Value = 1570.6086
$ws.Range("I107").Value = 855.375
$ws.Range("J107").Value = 1952.0667
$ws.Range("K107").Value = 855.375
$ws.Range("L107").Value = 1952.0667
$ws.Range("M107").Value = 1064.625
$ws.Range("N107").Value = -5792.0667
$ws.Range("H132").Value = 4597
$ws.Range("I132").Value = 8198.799999999999
$ws.Range("J132").Value = 3596.5
$ws.Range("K132").Value = 24596.4
$ws.Range("L132").Value = 10789.5
$ws.Range("M132").Value = -22066.4
$ws.Range("N132").Value = -15849.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2356.8667
$ws.Range("J61").Value = 7005
$ws.Range("L61").Value = 7005
$ws.Range("N61").Value = -7409
$ws.Range("H68").Value = 50001680
$ws.Range("I68").Value = 1465.1818
$ws.Range("J68").Value = 111113060
$ws.Range("K68").Value = 1465.1818
$ws.Range("L68").Value = 111113060
$ws.Range("M68").Value = -716.1818000000001
$ws.Range("N68").Value = -111114558
$ws.Range("H71").Value = 50001680
$ws.Range("I71").Value = 1465.1818
$ws.Range("J71").Value = 111113060
$ws.Range("K71").Value = 7325.909000000001
$ws.Range("L71").Value = 555565300
$ws.Range("M71").Value = -3581.909000000001
$ws.Range("N71").Value = -555572788
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 2356.8667
$ws.Range("J113").Value = 7005
$ws.Range("L113").Value = 7005
$ws.Range("N113").Value = -11345
$ws.Range("H132").Value = 22815188
$ws.Range("I132").Value = 48161816
$ws.Range("K132").Value = 144485448
$ws.Range("M132").Value = -144482918

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 29862.2
$ws.Range("J123").Value = 29862.2
$ws.Range("L123").Value = 29862.2
$ws.Range("N123").Value = -39662.2
$ws.Range("H132").Value = 2186.4783
$ws.Range("J132").Value = 2728.3572
$ws.Range("L132").Value = 8185.071599999999
$ws.Range("N132").Value = -13245.0716
$ws.Range("H136").Value = 1355.6522
$ws.Range("I136").Value = 764.1667
$ws.Range("K136").Value = 2292.5001
$ws.Range("M136").Value = 257.4998999999998
